# Actualizacion automatica hashcode vie mar  6 01:38:43 CET 2020
# Updates the hashcode column (column B) for the rows whose code (column A)
# matches one of the entries below, replacing the stale hash with the new one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 11; Code = "05-050301A"; Old = "10886ff983b31d2b23a61037bee03245"; New = "cf2f1f7db0e2df3878213475336921a3" },
    @{ Row = 17; Code = "05-050305TP"; Old = "038304534db84b6ed92e4778656351d4"; New = "0ac9d1f217e8aa92141b97ddd2e464a5" },
    @{ Row = 24; Code = "05-050316TC"; Old = "a555b1001cb4eb4774155e4e0d45ed42"; New = "9e002fac2a2d51c95fa2def6a29df35e" },
    @{ Row = 29; Code = "05-050302A"; Old = "cee7288738f2d18a531208d0aff3184a"; New = "020895d6d25807ceea798e4b8cd227b7" },
    @{ Row = 34; Code = "05-050316TP"; Old = "199a671aa92b1cddc0ce99fa4e4b1e95"; New = "85b945b689bab52ab952fc0bdd4056db" },
    @{ Row = 126; Code = "05-050309A"; Old = "1dcc289a9d0f869602d31804cf60a0ee"; New = "0e7449a6be04ef7efd69afaf0df094cc" },
    @{ Row = 133; Code = "05-050312TP"; Old = "79d1d0c0b7ca3ccbf7b1e9c227371f3f"; New = "a1f0aeac8802c0250624fc9fa2c26529" },
    @{ Row = 136; Code = "05-050312TC"; Old = "ced153bef9faf7a242b0bc254c1cbd1b"; New = "f2ed162b9f08ab469f340bd77d18c52a" },
    @{ Row = 159; Code = "05-050203TP"; Old = "3179b1019818ad8c556b64072a9463e9"; New = "aaac34bf96dd1a6394dd3ca31665d0c9" },
    @{ Row = 162; Code = "05-050308A"; Old = "fac192a900ed093137d7272371060418"; New = "f6e10bcb8d47e08bc2d03119866ebb46" },
    @{ Row = 175; Code = "05-050303TP"; Old = "193cd4c5a8bab3fca87960a3d4334401"; New = "1de777d94e43086bcd7f3ab707a163fe" },
    @{ Row = 180; Code = "05-050303TC"; Old = "5abe0996962ce49df8ad4ecad6d6e6b1"; New = "79f57dada9cae6290c1593b3d8025a1c" },
    @{ Row = 183; Code = "05-050305A"; Old = "1d5ad6fd7ba9d2853a1f6b91fecc2317"; New = "0ca4f74849b3b000caf79995a5889750" },
    @{ Row = 191; Code = "05-050314TP"; Old = "1083da5df02bf38f818a271508322574"; New = "32cfcd119d179ac4a5597dc259240032" },
    @{ Row = 198; Code = "05-050314TC"; Old = "8067240336eb47712eaecf0e3379c696"; New = "218495e0ce5c193e6de4326bb103aa11" },
    @{ Row = 200; Code = "05-050306A"; Old = "958cbe82ce1b73a5e002af2bbbfd2c6b"; New = "a84840917c81c5a306c6fab2c73dd40c" },
    @{ Row = 213; Code = "05-050303A"; Old = "3d3e8d23a97d243c3fb637cfccec89d7"; New = "8fe3858276d30b3ad2d6cc503d53cf9d" },
    @{ Row = 227; Code = "05-050205TP"; Old = "4f00912ce8da1abd59079cc7ea9c841a"; New = "fa6365624fd8a90e2d09bcf1285e731d" },
    @{ Row = 228; Code = "05-050304A"; Old = "4648fda5675865d389bb840e5e81c602"; New = "ac8bed145257d1de473f50ed7b7c69cc" },
    @{ Row = 232; Code = "05-050205TC"; Old = "07166dd52e9e8ba9e7603046b8a011db"; New = "772259b09e9211f2c703127cc02e4a77" },
    @{ Row = 281; Code = "05-050201TC"; Old = "882631270fdb637858541421070c4b7e"; New = "d47b4c2c37695aeaedf46052fc07213c" },
    @{ Row = 302; Code = "05-050310TP"; Old = "20ccd7545bd40ce35161e419bfdcfe56"; New = "e6e9d639920f8cd3c896404d721d0c33" },
    @{ Row = 339; Code = "05-050201TP"; Old = "addf5a747b264949fa9ae8e691ca5087"; New = "0cfcf0cdbc873d2da6b6d2d79315cafe" },
    @{ Row = 464; Code = "05-050204A"; Old = "b189b4d6b4454b07494170016cc0a052"; New = "878b36f81cb1d4e4d02d61d7f76f4b38" },
    @{ Row = 483; Code = "05-050205A"; Old = "f8d0b1e4bb52fb431e0dc3e7cf49c690"; New = "db0ef39e011ea89708a4b5d76f64a6bb" },
    @{ Row = 485; Code = "05-050314A"; Old = "3dab421690256830d891eb1dbd6545c8"; New = "137211fbc02800389c315e3667e2f3e3" },
    @{ Row = 507; Code = "05-050311A"; Old = "1cacb3cee02312b2a93c65a2a344c7bf"; New = "ec0daeee44dcb5f4d9d1e9e844866c45" },
    @{ Row = 513; Code = "05-050306TP"; Old = "57453290a028d0832d2d6a87aba3f3d1"; New = "412827c9b7c01576f8915108af0677e6" },
    @{ Row = 521; Code = "05-050317TC"; Old = "7bec1385342fed9aa75716535350b327"; New = "c4ef335f26aaebd49a24f54269269a74" },
    @{ Row = 522; Code = "05-050312A"; Old = "cfc370cdddbfda7f07d64e0347ab9971"; New = "7e44aa3932207675840f13940dfded53" },
    @{ Row = 532; Code = "05-050317TP"; Old = "bd765d93499de8a428406c20c7de6700"; New = "f21adc895e78d6e312e03aa4f3fd153e" },
    @{ Row = 558; Code = "05-050310A"; Old = "54851fd6c0d1a26d2ae9d06c37d1fbb8"; New = "ce6e8ce7390f23145e4e60636e1540db" },
    @{ Row = 624; Code = "05-050204TP"; Old = "0690257d524fa65e2c39a24884c7519e"; New = "bf03cfbb025a35e673c7755066dadfab" },
    @{ Row = 635; Code = "05-050204TC"; Old = "b984c87dcf8554dba12699230be4fd78"; New = "8ee5eefcebf6ffd22aa05877f91ffb19" },
    @{ Row = 637; Code = "05-050302TP"; Old = "4ffb3ea8d532b90ba41ae1b4caeab26c"; New = "93ef2328a3b5c2a9f75453d8c4ad9cbd" },
    @{ Row = 674; Code = "05-050317A"; Old = "5e8a3bfaf7d985c7619ed91006c40ba5"; New = "bc59dbcf49c7b976c956c9f08f6a5d95" },
    @{ Row = 708; Code = "05-050304TC"; Old = "8c495fbb95df2dafbecb6aa9fd1f317a"; New = "ac11ff5172c43564a5b15233fd7c3275" },
    @{ Row = 712; Code = "05-050315A"; Old = "eed96e3ef1c25fb650d56cd4b8d8dc26"; New = "ae5edf4d2f601d2c7a0da70f96a2044b" },
    @{ Row = 723; Code = "05-050304TP"; Old = "82fed5dc3184faa31d69a830a10ada13"; New = "4ddd244a02ae194577a8d7a8096c1357" },
    @{ Row = 734; Code = "05-050315TC"; Old = "9d2ff17e621f36576bd240e35116c805"; New = "c669b0de5ff0c1ea604915c1706f4995" },
    @{ Row = 737; Code = "05-050316A"; Old = "7ac51dffd4b9e5f46303f624a41708a5"; New = "72e1dffbd0d0ec525203283740ded950" },
    @{ Row = 750; Code = "05-050315TP"; Old = "2958a31aa257329ad526cdbdc3e9be0e"; New = "b4f3e60af1ba9d42c64f7f9829fd7475" },
    @{ Row = 838; Code = "05-050311TC"; Old = "609e13c97c8ea9422fcd925b50c0bb4f"; New = "1bbbfb0d19de04d5f54fcdff7ffbcf96" },
    @{ Row = 843; Code = "05-050311TP"; Old = "2ee6460c61db675a0c438b7cc8ca8745"; New = "7be205ce29f91ce3180a07bd6c517373" },
    @{ Row = 862; Code = "05-050309TC"; Old = "2549441feec73fad726ef2286fad0e82"; New = "cd493707b130401d49743f9fb5054dd6" },
    @{ Row = 870; Code = "05-050309TP"; Old = "2ac976d0abbdb6753b1e8028cc220b23"; New = "2868f8250a17e53d0e7b5226a008fd5f" }
)

foreach ($u in $updates) {
    $codeCell = $ws.Cells.Item($u.Row, 1)
    $hashCell = $ws.Cells.Item($u.Row, 2)
    if ($codeCell.Value() -eq $u.Code -and $hashCell.Value() -eq $u.Old) {
        $hashCell.Value = $u.New
    }
}
